$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Number of trials"), shifting
# everything from C onward one column to the right.
$ws.Columns.Item(3).Insert()

# Header for the new column.
$ws.Range("C1").Value = "Resting Rate"

# Resting Rate values per row.
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 1.7
$ws.Range("C4").Value = 6.03
$ws.Range("C5").Value = 7.05
$ws.Range("C6").Value = 0

# Restore the selection state to match the saved workbook.
$ws.Range("E17").Select()
